# Applies the "final draft files to be submitted" edit to the Crowdfunding
# Discussion document.
#
# wdReplace constants used below:
#   wdFindContinue = 1   (Wrap parameter to Find.Execute)
#   wdReplaceAll   = 2   (Replace parameter to Find.Execute)

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $find"
    }
    return $result
}

# --- Bullet: "While the theater parent category..." -----------------------
Replace-Text 'of successful campaigns, the “journalism”' 'of successful campaigns (187 of 344), the “journalism”'
Replace-Text 'highest rate of successful campaigns; granted' 'highest rate of successful campaigns (4 of 4); granted'

# --- Bullet: "Campaigns started in June or July..." ------------------------
Replace-Text 'more likely to be successful, with 11.3% of the campaigns (113 campaigns)' 'more likely to be successful (63.2%), with 11.3% of the total campaigns (113 campaigns)'

# --- Bullet: "I would want the individual donations table..." --------------
Replace-Text 'more likely to receive donations need the beginning' 'more likely to receive donations near the beginning'

# --- Bullet: "I would also investigate the staff_pic..." -------------------
# Text is unchanged; only the proofing-error markers that bracket "staff_pic"
# are removed. Running it through Find/Replace (even as a no-op textual
# change) collapses the split runs/proofErr markers into a single plain run,
# matching the target.
Replace-Text 'I would also investigate the staff_pic and spotlight' 'I would also investigate the staff_pic and spotlight'

# --- Paragraph: "For both the successful and failed outcomes..." + new ----
# bullet-point paragraphs appended right after it (with a blank paragraph
# in between), all applied atomically so run merging behaves predictably.
Replace-Text 'For both the successful and failed outcomes vs. backers_count, the median is a better measure because the data is right skewed. As you can see in the graphs I created in the “Statistical Analysis” tab in the workbook, with most of the data living in the first bucket on the left with many outliers (~8% of the data in both cases). As these outliers are very large compared to the rest of the data, the mean is larger than the median, and therefore less representative of the data.' 'For both the successful and failed outcomes vs. backers_count, the median is a better measure because the data is right skewed. As you can see in the graphs I created in the “Statistical Analysis” tab in the workbook, most of the data lives in the first bucket on the left with many outliers (~8% of the data in both cases are outliers). As these outliers are very large compared to the rest of the data, the mean is larger than the median, and therefore less representative of the data.^p^pThere is more variability within the successful campaign data, because the standard deviation, variance, range, and interquartile range are all higher than their counterparts within the failed campaign data. This makes sense as for a campaign to fail, the percentage funded must be somewhere between 0% and 100%, inclusive. However, a campaign that succeeds is one that is funded at 100% or more, and there are many highly funded. For example, there are 24 campaigns that were funded at 1000% percent or more within the data, with another 66 campaigns funded between 500% and 1000%.'

# --- Style: DefaultParagraphFont becomes semi-hidden -----------------------
$style = $d.Styles("Default Paragraph Font")
$style.SemiHidden = $true

Write-Output "done"
